# Append the 2025-04-22 price row (row 52) to each of the 9 Solar Prices
# sheets, matching the existing "Date"/"Price" text-column layout.

$wb = $excel.ActiveWorkbook

$newDate = "2025-04-22"

$rows = @(
    @{ Sheet = "N-Dense";                    Price = "39.5" },
    @{ Sheet = "N-Type";                     Price = "39" },
    @{ Sheet = "N-type Wafer";               Price = "1.23" },
    @{ Sheet = "Cell Topcon 183mm";          Price = "0.293" },
    @{ Sheet = "Module Topcon 183mm";        Price = "0.09" },
    @{ Sheet = "Silver Rear_side";           Price = "5,367" },
    @{ Sheet = "Silver Busbar front-side";   Price = "8,035" },
    @{ Sheet = "Silver finger front-side";   Price = "8,085" },
    @{ Sheet = "USD_CNY";                    Price = "7.3133" }
)

foreach ($row in $rows) {
    $ws = $wb.Worksheets.Item($row.Sheet)

    $target = $ws.Range("A52:B52")

    # Force the new cells to be stored as plain text (matching every other
    # cell in these columns) instead of letting Excel auto-convert the
    # date-looking / number-looking strings into a real date serial or
    # number, then strip the temporary "Text" number format back off so no
    # stray formatting is left behind on the new row.
    $target.NumberFormat = "@"
    $ws.Range("A52").Value = $newDate
    $ws.Range("B52").Value = $row.Price
    $target.ClearFormats()
}

Write-Output "Appended 2025-04-22 row to $($rows.Count) sheets"
